$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New shared strings must be created in this order so the sharedStrings
# table grows as: ... 36 (existing), 37 "Não procede", 38 "Analisado" ---
$ws.Range("E18").Value = "Não procede"
$ws.Range("E26").Value = "Não procede"

$ws.Range("C17").Value = "Desenvolvido"
$ws.Range("C18").Value = "Analisado"
$ws.Range("C26").Value = "Analisado"

# --- Re-colour rows 17, 18 and 26 from the "orange" style to the "green"
# style already used elsewhere on the sheet (e.g. row 19), by copying the
# cell format from a green cell so the existing style index is reused
# instead of a brand-new one being fabricated. ---
$ws.Range("A19:D19").Copy()
$ws.Range("A17:D17").PasteSpecial(-4122)
$ws.Range("A18:D18").PasteSpecial(-4122)
$ws.Range("A26:D26").PasteSpecial(-4122)

$ws.Range("D19").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("E26").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Restore the view: scroll back to the top and select A22 ---
$ws.Range("A22").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
